# Update the date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-01-31 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-02-01 Thursday", 2)

# Update the multiplication-answer table. Cells are addressed by
# (row, column) so that duplicate answer values (e.g. 527x3=1581 and
# 272x4=1088 each show up both as an old and a new value) are never
# confused with each other via text search.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "742×7=5194" },
    @{ Row = 1;  Col = 2; Text = "527×3=1581" },
    @{ Row = 1;  Col = 3; Text = "845×9=7605" },
    @{ Row = 1;  Col = 4; Text = "246×7=1722" },
    @{ Row = 1;  Col = 5; Text = "563×6=3378" },

    @{ Row = 5;  Col = 1; Text = "765×3=2295" },
    @{ Row = 5;  Col = 2; Text = "138×2=276" },
    @{ Row = 5;  Col = 3; Text = "512×6=3072" },
    @{ Row = 5;  Col = 4; Text = "826×5=4130" },
    @{ Row = 5;  Col = 5; Text = "813×5=4065" },

    @{ Row = 10; Col = 1; Text = "461×8=3688" },
    @{ Row = 10; Col = 2; Text = "425×3=1275" },
    @{ Row = 10; Col = 3; Text = "517×2=1034" },
    @{ Row = 10; Col = 4; Text = "797×4=3188" },
    @{ Row = 10; Col = 5; Text = "479×3=1437" },

    @{ Row = 15; Col = 1; Text = "244×9=2196" },
    @{ Row = 15; Col = 2; Text = "890×9=8010" },
    @{ Row = 15; Col = 3; Text = "457×3=1371" },
    @{ Row = 15; Col = 4; Text = "940×5=4700" },
    @{ Row = 15; Col = 5; Text = "701×4=2804" },

    @{ Row = 20; Col = 1; Text = "888×6=5328" },
    @{ Row = 20; Col = 2; Text = "905×6=5430" },
    @{ Row = 20; Col = 3; Text = "379×2=758" },
    @{ Row = 20; Col = 4; Text = "272×4=1088" },
    @{ Row = 20; Col = 5; Text = "110×8=880" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Output "done"
